$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F. This shifts the old F,G,H columns to G,H,I
# and automatically updates the dataValidations sqref ranges.
$ws.Columns("F").Insert()

# --- New column F holds a shared formula that reproduces the expense total ---
# F1 is the master formula cell (the shared-formula anchor covers F1:F2).
$ws.Range("F1:F2").Formula = '=if(And(G1<>"",H1<>""),if(E1<>"",E1,D1),)'

# F1 styling: same font/fill as the adjoining validated cell (I1), but with a
# right-aligned "#,##0.00" number format.
$ws.Range("I1").Copy()
$ws.Range("F1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F1").NumberFormat = "#,##0.00"
$ws.Range("F1").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

# F2 styling: same base look, but "#,##0.00;(#,##0.00)" number format.
$ws.Range("I2").Copy()
$ws.Range("F2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F2").NumberFormat = "#,##0.00;(#,##0.00)"
$ws.Range("F2").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight
